$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45310, now 45311)
$ws.Range("A1").Value = 45311

# Update price column D for rows 34-37 ("CON TOPE" PICC items)
$ws.Range("D34").Value = 1304.725
$ws.Range("D35").Value = 1739.632
$ws.Range("D36").Value = 2519.906
$ws.Range("D37").Value = 2143.067

# Update price column D for rows 41-44 ("SIN TOPE" PICS items)
$ws.Range("D41").Value = 1675.673
$ws.Range("D42").Value = 2213.895
$ws.Range("D43").Value = 2788.521
$ws.Range("D44").Value = 2437.246

# Re-create merged ranges so the A1:G1 merge ends up recorded before the
# A39:D39 merge (matches the refreshed merge order after re-touching A1).
$ws.Range("A30:F30").UnMerge()
$ws.Range("A32:D32").UnMerge()
$ws.Range("A39:D39").UnMerge()
$ws.Range("A1:G1").UnMerge()

$ws.Range("A30:F30").Merge()
$ws.Range("A32:D32").Merge()
$ws.Range("A1:G1").Merge()
$ws.Range("A39:D39").Merge()
